$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting existing rows 106-132 down to 107-133.
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with a new weekly record for Mango.
$ws.Cells.Item(106, 1).Value = 5
$ws.Cells.Item(106, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(106, 3).Value = "Maule"
$ws.Cells.Item(106, 4).Value = 44736
$ws.Cells.Item(106, 5).Value = 7
$ws.Cells.Item(106, 6).Value = "Fruta"
$ws.Cells.Item(106, 7).Value = 100108
$ws.Cells.Item(106, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(106, 9).Value = 100108002
$ws.Cells.Item(106, 10).Value = "Mango"
$ws.Cells.Item(106, 11).Value = "Sin especificar"
$ws.Cells.Item(106, 12).Value = "Primera"
$ws.Cells.Item(106, 13).Value = 242
$ws.Cells.Item(106, 14).Value = 8000
$ws.Cells.Item(106, 15).Value = 8000
$ws.Cells.Item(106, 16).Value = 8000
$ws.Cells.Item(106, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(106, 18).Value = "Brasil"
$ws.Cells.Item(106, 19).Value = 2000
$ws.Cells.Item(106, 20).Value = 4
